# Performance average for 256 simulations workflow
$wb = $excel.ActiveWorkbook

# --- TTC sheet: update the "256 simulations" (column M/N/O/P @ rows 2,4-11)
#     throughput figures from the old per-run numbers to a uniform 870,
#     and repoint the H-column percentage formulas at column M instead of O.
$ttc = $wb.Worksheets.Item("TTC")

foreach ($r in 2,4,5,6,7,8,9,10,11) {
    $ttc.Range("M$r`:P$r").Value = 870
}

$ttc.Range("H4").Formula = "=(M4/B4)*100"
$ttc.Range("H5").Formula = "=(M5/B5)*100"
$ttc.Range("H6").Formula = "=(M6/B6)*100"
$ttc.Range("H7").Formula = "=(M7/B7)*100"

# H11 becomes part of the formula's fill range but stays blank (styled only)
$ttc.Range("H11").NumberFormat = "0"

# Make TTC the active sheet/tab, with H4 selected (matches the saved view)
$ttc.Activate()
$ttc.Range("H4").Select()

# --- Tw / Te sheets: selection changes to a block A1:B7 (no explicit active cell)
$tw = $wb.Worksheets.Item("Tw")
$tw.Range("A1:B7").Select()

$te = $wb.Worksheets.Item("Te")
$te.Range("A1:B7").Select()

# Restore TTC as the active sheet/tab (it must be the one left selected)
$ttc.Activate()
$ttc.Range("H4").Select()
